$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 21 - "Combatsystem for Axe" already has a task name; fill in estimate/actual hours and responsible person
$ws.Range("B21").Value = 6
$ws.Range("C21").Value = 6
$ws.Range("F21").Value = "Cedric"

# Row 23 - new task "Sprite für Springen"
$ws.Range("A23").Value = "Sprite für Springen"
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = "Cedric"

# Row 24 - new task "Ax in the Stone for Cave"
$ws.Range("A24").Value = "Ax in the Stone for Cave"

# Update selection to match the commit
$ws.Range("A24").Select()
